$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Copies served by Dep. Clerk ___ ... ___ to:" -- lengthen the two blank
#    lines (underscores) used for the deputy clerk's initials and the date.
# ---------------------------------------------------------------------------
$old1 = "Copies served by Dep. Clerk ___________ on the following date ___________ to:"
$new1 = "Copies served by Dep. Clerk ___________________________ on the following date ____________________ to:"
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "Prosecutor's Office: PS     OM ..." -- drop one space from the first
#    "PS___OM" gap only (the other two "PS___OM" / "PS___EM" gaps on the
#    same line are untouched).
# ---------------------------------------------------------------------------
$old2 = "Office: PS     OM"
$new2 = "Office: PS    OM"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Community Control / County Jail service line.
#    Originally two runs split by a <w:tab/>:
#       Run A: "{% if community_control... %}Community Control: PS "
#       Run B: <tab/>"EM;{% endif %}"
#    Target: a single run whose text is Run A (unchanged) + the tab turned
#    into 3 literal spaces (4 spaces total before "EM;") + "EM;{% endif %}"
#    + a brand-new "{% if jail_terms... %}County Jail: PS   EM;{% endif %}"
#    clause appended right after.
#    We locate the start of "Community Control: PS" and the end of
#    "EM;{% endif %}" with two independent Find operations, build the Range
#    spanning both (which also swallows the intervening tab character), and
#    overwrite its .Text in one shot -- Word naturally collapses that into a
#    single run using the formatting already present.
# ---------------------------------------------------------------------------
$partStart = $d.Content
$foundStart = $partStart.Find.Execute("Community Control: PS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$partEnd = $d.Content
$foundEnd = $partEnd.Find.Execute("EM;{% endif %}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($foundStart -and $foundEnd) {
    $span = $d.Range($partStart.Start, $partEnd.End)
    $span.Text = "Community Control: PS    EM;{% endif %}{% if jail_terms.ordered is true or apply_jtc == ‘Sentence’ %}County Jail: PS   EM;{% endif %}"
}

# ---------------------------------------------------------------------------
# 4. styles.xml <w:latentStyles><w:lsdException w:name="Table Grid" .../>
#    should lose its semiHidden="1"/unhideWhenUsed="1" flags (left with only
#    uiPriority="59"). This per-document "latent style" registry isn't part
#    of the Word object model that real Word/VBA (or this host) exposes, so
#    there is no supported COM call that can reach it. Guard defensively in
#    case a future host revision adds support, without failing the script.
# ---------------------------------------------------------------------------
try {
    $latent = $d.LatentStyles
    if ($latent -ne $null) {
        $tableGrid = $latent.Item("Table Grid")
        if ($tableGrid -ne $null) {
            $tableGrid.SemiHidden = $false
            $tableGrid.UnhideWhenUsed = $false
        }
    }
} catch {
    # Not supported by the object model exposed here; nothing further to do.
}

Write-Output "Edits applied."
